$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.843.58"
$ws.Range("E2").Value = "  -1.07%  "
$ws.Range("D3").Value = "1.896.20"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7782"
$ws.Range("E5").Value = "  +5.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.02"
$ws.Range("E6").Value = "  -1.87%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3066"
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.52"
$ws.Range("E9").Value = "  -5.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06859"
$ws.Range("E10").Value = "  -2.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07985"
$ws.Range("D12").Value = "1.916.85"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7367"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.175"
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.38"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").Value = "29.851.85"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.78"
$ws.Range("E17").Value = "  -4.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.886"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.10"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007696"
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "2.153.00"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.929"
$ws.Range("E24").Value = "  -3.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.73"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.273"
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1318"
$ws.Range("E27").Value = "  +2.21%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.74"
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.026"
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.390"
$ws.Range("E30").Value = "  +2.60%  "
$ws.Range("E31").Value = "  -2.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.269"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.070"
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05261"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("E35").Value = "  -4.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7283"
$ws.Range("E36").Value = "  -3.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.725"
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01905"
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.174"
$ws.Range("E40").Value = "  -3.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4414"
$ws.Range("E41").Value = "  -2.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.03"
$ws.Range("E42").Value = "  -4.10%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8372"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.874"
$ws.Range("E45").Value = "  -4.69%  "
$ws.Range("E46").Value = "  -3.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.28"
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.742"
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("D49").Value = "2.060.64"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.08"
$ws.Range("E50").Value = "  -3.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "924.01"
$ws.Range("E51").Value = "  -1.35%  "
